# Update "want to go" counts (column F) across the four sheets to the
# refreshed snapshot values captured at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 2312
$ws.Range("F3").Value = 386
$ws.Range("F4").Value = 195
$ws.Range("F5").Value = 307
$ws.Range("F6").Value = 307
$ws.Range("F10").Value = 529
$ws.Range("F11").Value = 727
$ws.Range("F14").Value = 379
$ws.Range("F16").Value = 996
$ws.Range("F17").Value = 17780
$ws.Range("F18").Value = 408
$ws.Range("F20").Value = 185
$ws.Range("F22").Value = 161
$ws.Range("F25").Value = 158
$ws.Range("F26").Value = 15
$ws.Range("F27").Value = 305
$ws.Range("F28").Value = 121

# --- Sheet "演出" (Performance) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 167
$ws.Range("F6").Value = 195
$ws.Range("F8").Value = 3323
$ws.Range("F10").Value = 37
$ws.Range("F16").Value = 2729

# --- Sheet "本地生活" (Local Life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 73
$ws.Range("F4").Value = 517

# --- Sheet "全部类型" (All Types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 73
$ws.Range("F6").Value = 2312
$ws.Range("F7").Value = 517
$ws.Range("F8").Value = 386
$ws.Range("F9").Value = 195
$ws.Range("F10").Value = 307
$ws.Range("F11").Value = 307
$ws.Range("F13").Value = 167
$ws.Range("F16").Value = 195
$ws.Range("F19").Value = 529
$ws.Range("F20").Value = 727
$ws.Range("F23").Value = 379
$ws.Range("F25").Value = 996
$ws.Range("F26").Value = 17780
$ws.Range("F28").Value = 3323
$ws.Range("F30").Value = 37
$ws.Range("F32").Value = 408
$ws.Range("F34").Value = 185
$ws.Range("F38").Value = 161
$ws.Range("F43").Value = 159
$ws.Range("F44").Value = 15
$ws.Range("F45").Value = 305
$ws.Range("F46").Value = 121
$ws.Range("F47").Value = 2729
